$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numeric-looking values (e.g. "20.01", "0.550")
# as plain text in the source data, not as numbers. Assigning such a string
# straight to .Value lets Excel auto-coerce it into a real number (losing
# trailing zeros / introducing floating-point noise, e.g. "0.550" -> 0.55,
# "20.01" -> 20.010000000000002). Pre-formatting each target cell as Text
# keeps the write literal; formatting is cleared again afterwards so no
# extra number-format/quote-prefix styling lingers on the cell.
$textCells = @("D2", "D3", "D5", "D10", "D11", "D13", "D14", "D15", "D17", "D18", "D20", "D21", "D23", "D24", "D28", "D29", "D31", "D34", "D37", "D38", "D40", "D42", "D45", "D46", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.242.50"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "1.645.71"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "217.19"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D10").Value = "20.01"
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.873.13"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "1.620.76"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "0.550"
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "63.55"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "26.219.74"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "195.89"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").Value = "4.44"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").Value = "6.35"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "143.65"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").Value = "6.96"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "15.62"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").Value = "0.0506"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D34").Value = "1.61"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.135.19"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.554"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("D40").Value = "0.0158"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "5.67"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").Value = "1.782.49"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").Value = "56.35"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("E47").Value = "  +3.90%  "
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("D49").Value = "7.72"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  +1.26%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}

